$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = "W1-D1(24/1)"
$ws.Cells.Item(1, 3).Value = "W1-D2(26/1)"
$ws.Cells.Item(1, 4).Value = "W1-D3(28/1)"
$ws.Cells.Item(1, 5).Value = "W1-D4(30/1)"
$ws.Cells.Item(1, 6).Value = "W1-D5(31/1)"
$ws.Cells.Item(1, 7).Value = "W1-D6(3/2)"
$ws.Cells.Item(1, 8).Value = "W1-D7(4/2)"
$ws.Cells.Item(1, 9).Value = "W2-D1(6/2)"
$ws.Cells.Item(1, 10).Value = "W2-D2(7/2)"
$ws.Cells.Item(1, 11).Value = "W2-D3(11/2)"
$ws.Cells.Item(1, 12).Value = "W2-D4(12/2)"
$ws.Cells.Item(1, 13).Value = "W2-D5(13/2)"
$ws.Cells.Item(1, 14).Value = "W2-D6(14/2)"
$ws.Cells.Item(1, 15).Value = "W2-D7(16/2)"
$ws.Cells.Item(1, 16).Value = "W3-D1(17/2)"
$ws.Cells.Item(1, 17).Value = "W3-D2(18/2)"
$ws.Cells.Item(1, 18).Value = "W3-D3(19/2)"
$ws.Cells.Item(1, 19).Value = "W3-D4(20/2)"
$ws.Cells.Item(1, 20).Value = "W3-D5(21/2)"
$ws.Cells.Item(1, 21).Value = "W3-D6(23/2)"
$ws.Cells.Item(1, 22).Value = "W3-D7(24/2)"
$ws.Cells.Item(1, 23).Value = "W4-D1(25/2)"
$ws.Cells.Item(1, 24).Value = "W4-D2(26/2)"
$ws.Cells.Item(1, 25).Value = "W4-D3(27/2)"
$ws.Cells.Item(1, 26).Value = "W4-D4(28/2)"
$ws.Cells.Item(1, 27).Value = "W4-D5(2/3)"
$ws.Cells.Item(1, 28).Value = "W4-D6(3/3)"
$ws.Cells.Item(1, 29).Value = "W4-D7(4/3)"
$ws.Cells.Item(1, 30).Value = "W5-D1(5/3)"
$ws.Cells.Item(1, 31).Value = "W5-D2(6/3)"
$ws.Cells.Item(1, 32).Value = "W5-D3(7/3)"
$ws.Cells.Item(1, 33).Value = "W5-D4(9/3)"
$ws.Cells.Item(1, 34).Value = "W5-D5(10/3)"
$ws.Cells.Item(1, 35).Value = "W5-D6(11/3)"
$ws.Cells.Item(1, 36).Value = "W5-D7(12/3)"
$ws.Cells.Item(1, 37).Value = "W6-D1(13/3)"
$ws.Cells.Item(1, 38).Value = "W6-D2(14/3)"
$ws.Cells.Item(1, 39).Value = "W6-D3(16/3)"
$ws.Cells.Item(1, 40).Value = "W6-D4(17/3)"
$ws.Cells.Item(1, 41).Value = "W6-D5(18/3)"
$ws.Cells.Item(1, 42).Value = "W6-D6(19/3)"
$ws.Cells.Item(1, 43).Value = "W6-D7(20/3)"
